# export-structure-stats.xlsx - add two new rows ("Orientation vers CIAS" and
# "Autre orientation") to the "Répartition des orientations" block on the
# "Stats structure" sheet, just below the existing "Orientation vers
# Organisme agrée" row, pushing the "3. TOTAL DES INTERACTIONS ..." block
# (and everything below it) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stats structure")

# Row 102 already exists (blank, but pre-formatted). Insert one brand-new
# row at 103 so the "3. TOTAL DES INTERACTIONS ..." section (old row 104)
# and everything after it shifts down to make room for the two new
# orientation rows (102 and 103).
$ws.Rows("103:103").Insert()

# Fill in the two new orientation rows.
$ws.Range("B102").Value = "Orientation vers CIAS"
$ws.Range("B103").Value = "Autre orientation"

# Match the row height used by the surrounding rows in this block.
$ws.Rows("102:102").RowHeight = 16
$ws.Rows("103:103").RowHeight = 16

# Update the active selection / scroll position to match where the edit
# was made.
$ws.Range("B101").Select()
